# Update schedule: Monday/Tuesday/Wednesday shifts now end at 16:00 (8 hours)
# instead of 17:00 (9 hours). Affects rows 4,5,6 / 11,12,13 / 18,19,20 /
# 25,26,27 / 32 (the Mon/Tue/Wed rows for each week of November 2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 5, 6, 11, 12, 13, 18, 19, 20, 25, 26, 27, 32)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "08:00 - 16:00"
    $ws.Range("E$r").Value = 8
}

# Recalculate totals for the affected rows (13 rows each losing 1 hour).
$ws.Range("E33").Value = 257
$ws.Range("E34").Value = 4626
